# Atualização automática de CANELA.xlsx
# - Rename "Paineis DARQ" -> "PAINEIS DARQ"
# - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove the "Desarquivamentos Pendentes" sheet
# - Keep "PAINEIS DARQ" as the active sheet

$wb = $excel.ActiveWorkbook

# Suppress the "permanently delete" confirmation dialog for the sheet removal.
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
$null = $wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true

# Restore the original active/selected sheet.
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
